$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Item List")
$ws.Range("K1").Value = "X"
$ws.Range("K1").Font.Bold = $true
$ws.Range("K1").Borders(10).Weight = -4138
$ws.Range("K1").Borders(10).Color = 13421772
